$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 20 with the new class data (Height, Weight, Birthplace, #Siblings)
$ws.Range("A20").Value = 74
$ws.Range("B20").Value = 190
$ws.Range("C20").Value = "MD"
$ws.Range("D20").Value = 2

# Move the active selection to A21, matching the author's next data-entry position
$ws.Range("A21").Select()
